$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Correzione data rilascio revisione 00 (D2: 45656 -> 45626) ---
$ws.Range("D2").Value = 45626

# --- Add new table row for version 0.1 ---
$lr = $tbl.ListRows.Add()

# Copy the formatting of row 2 onto the new row 3 so per-cell styles
# (date format on C/D, wrap-text on B) line up with the rest of the table.
$ws.Range("A2:D2").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# A3 ("0.1") must be stored as text, matching how "0.0" is stored in A2.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "0.1"
$ws.Range("A3").ClearFormats()

# B3: Modifiche
$ws.Range("B3").Value = "Inserimento versione progetto"

# C3: Data inizio modifiche (16/12/2024)
$ws.Range("C3").Value = 45642

# D3: Data rilascio - not released yet, leave empty
$ws.Range("D3").Clear()

# --- UI selection state ---
[void]$ws.Range("K2").Select()
